$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: bump the generation Date
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2025-07-29T07:08:53+00:00"

# ------------------------------------------------------------------
# 2. Elements sheet: insert a new row for
#    FonctionQualifiee.exerciceProfessionnel right above the existing
#    FonctionQualifiee.fonctionQualifiee row (old row 6), pushing the
#    latter down to row 7.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Shift old row 6 down to row 7, inheriting its formatting.
$ws.Rows.Item(6).Insert()

# Give the freshly inserted row the same formatting (borders / wrap /
# vertical alignment) as the rows around it.
$ws.Range("A5:AJ5").Copy()
$ws.Range("A6:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content (FonctionQualifiee.exerciceProfessionnel).
$ws.Range("A6").Value2 = "FonctionQualifiee.exerciceProfessionnel"
$ws.Range("B6").Value2 = "FonctionQualifiee.exerciceProfessionnel"
$ws.Range("D6").Value2 = ""
$ws.Range("F6").Value2 = "1"
$ws.Range("G6").Value2 = "1"
$ws.Range("H6").Value2 = ""
$ws.Range("I6").Value2 = ""
$ws.Range("J6").Value2 = ""
$ws.Range("K6").Value2 = "Reference(https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/ExerciceProfessionnel)`n"
$ws.Range("L6").Value2 = "Lien vers la classe ExerciceProfessionnel."
$ws.Range("M6").Value2 = "Lien vers la classe ExerciceProfessionnel."
$ws.Range("P6").Value2 = ""
$ws.Range("R6").Value2 = ""
$ws.Range("S6").Value2 = ""
$ws.Range("T6").Value2 = ""
$ws.Range("U6").Value2 = ""
$ws.Range("V6").Value2 = ""
$ws.Range("W6").Value2 = ""
$ws.Range("X6").Value2 = ""
$ws.Range("Y6").Value2 = ""
$ws.Range("Z6").Value2 = ""
$ws.Range("AA6").Value2 = ""
$ws.Range("AB6").Value2 = ""
$ws.Range("AC6").Value2 = ""
$ws.Range("AD6").Value2 = ""
$ws.Range("AE6").Value2 = ""
$ws.Range("AF6").Value2 = "SavoirFaire.exerciceProfessionnel"
$ws.Range("AG6").Value2 = "1"
$ws.Range("AH6").Value2 = "1"
$ws.Range("AI6").Value2 = ""
$ws.Range("AJ6").Value2 = ""

# ------------------------------------------------------------------
# 3. Column widths: A/B and K grow to fit the new, longer content
#    (mirrors Excel's "best fit" recompute after the edit).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 31.5
$ws.Columns.Item(2).ColumnWidth = 31.5
$ws.Columns.Item(11).ColumnWidth = 74.666666666666667
